$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6757
$ws1.Range("F3").Value = 812
$ws1.Range("F5").Value = 139
$ws1.Range("F6").Value = 6
$ws1.Range("F7").Value = 716
$ws1.Range("F8").Value = 716
$ws1.Range("F9").Value = 12
$ws1.Range("F12").Value = 1104
$ws1.Range("F14").Value = 11
$ws1.Range("F15").Value = 701
$ws1.Range("F16").Value = 1014
$ws1.Range("F17").Value = 1340
$ws1.Range("F21").Value = 2
$ws1.Range("F22").Value = 556
$ws1.Range("F25").Value = 366
$ws1.Range("F27").Value = 1491
$ws1.Range("F29").Value = 528
$ws1.Range("F30").Value = 457
$ws1.Range("F31").Value = 450
$ws1.Range("F33").Value = 699
$ws1.Range("F34").Value = 1128
$ws1.Range("F36").Value = 2352
$ws1.Range("F38").Value = 1237
$ws1.Range("F41").Value = 3856
$ws1.Range("I36").Value = "//i0.hdslb.com/bfs/openplatform/202403/C8G3AOLM1709870001354.jpeg"

# --- Sheet: 演出 (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F19").Value = 4126
$ws2.Range("F24").Value = 2
$ws2.Range("F25").Value = 234

# --- Sheet: 本地生活 (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F5").Value = 1632
$ws3.Range("F8").Value = 959

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1632
$ws4.Range("F7").Value = 959
$ws4.Range("F8").Value = 6757
$ws4.Range("F10").Value = 812
$ws4.Range("F12").Value = 139
$ws4.Range("F13").Value = 6
$ws4.Range("F14").Value = 716
$ws4.Range("F15").Value = 716
$ws4.Range("F18").Value = 1104
$ws4.Range("F21").Value = 701
$ws4.Range("F25").Value = 1014
$ws4.Range("F26").Value = 1340
$ws4.Range("F30").Value = 556
$ws4.Range("F33").Value = 366
$ws4.Range("F35").Value = 1491
$ws4.Range("F37").Value = 528
$ws4.Range("F38").Value = 457
$ws4.Range("F39").Value = 450
$ws4.Range("F42").Value = 701
$ws4.Range("F43").Value = 1128
$ws4.Range("F45").Value = 2352
$ws4.Range("F49").Value = 1237
$ws4.Range("F51").Value = 3856
$ws4.Range("I45").Value = "//i0.hdslb.com/bfs/openplatform/202403/C8G3AOLM1709870001354.jpeg"
